$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15 (ALC)
$ws.Range("H15").Value = 186.65657
$ws.Range("I15").Value = 186.65657
$ws.Range("K15").Value = 559.96971
$ws.Range("M15").Value = -390.96971

# Row 76 (ALC)
$ws.Range("H76").Value = 4048.75
$ws.Range("I76").Value = 3596.6667
$ws.Range("J76").Value = 4320
$ws.Range("K76").Value = 3596.6667
$ws.Range("L76").Value = 4320
$ws.Range("M76").Value = -3281.6667
$ws.Range("N76").Value = -4950

# Row 79 (ALC)
$ws.Range("H79").Value = 4048.75
$ws.Range("I79").Value = 3596.6667
$ws.Range("J79").Value = 4320
$ws.Range("K79").Value = 3596.6667
$ws.Range("L79").Value = 4320
$ws.Range("M79").Value = -2504.6667
$ws.Range("N79").Value = -6504

# Row 87 (ALC)
$ws.Range("H87").Value = 99999.5
$ws.Range("J87").Value = 99999.5
$ws.Range("L87").Value = 99999.5
$ws.Range("N87").Value = -102495.5

# Row 90 (ALC)
$ws.Range("H90").Value = 99999.5
$ws.Range("J90").Value = 99999.5
$ws.Range("L90").Value = 299998.5
$ws.Range("N90").Value = -312478.5

# Row 107 (ALC)
$ws.Range("H107").Value = 84023.914
$ws.Range("I107").Value = 100648.8
$ws.Range("J107").Value = 899.5
$ws.Range("K107").Value = 100648.8
$ws.Range("L107").Value = 899.5
$ws.Range("M107").Value = -98728.8
$ws.Range("N107").Value = -4739.5

# Row 112 (ALC)
$ws.Range("H112").Value = 6513.185
$ws.Range("J112").Value = 7219
$ws.Range("L112").Value = 21657
$ws.Range("N112").Value = -23873

# Row 132 (ALC)
$ws.Range("H132").Value = 1443.8182
$ws.Range("I132").Value = 1313.2113
$ws.Range("J132").Value = 2989.3333
$ws.Range("K132").Value = 3939.6339
$ws.Range("L132").Value = 8967.999899999999
$ws.Range("M132").Value = -1409.6339
$ws.Range("N132").Value = -14027.9999

# Row 138 (ALC)
$ws.Range("H138").Value = 2672825.5
$ws.Range("I138").Value = 7411652
$ws.Range("J138").Value = 7235.854
$ws.Range("K138").Value = 22234956
$ws.Range("L138").Value = 21707.562
$ws.Range("M138").Value = -22229816
$ws.Range("N138").Value = -31987.562

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 55678.48
$ws.Range("I32").Value = 43521.383
$ws.Range("J32").Value = 119503.25
$ws.Range("K32").Value = 43521.383
$ws.Range("L32").Value = 119503.25
$ws.Range("M32").Value = -43234.383
$ws.Range("N32").Value = -120077.25

# Row 132 (ARM)
$ws.Range("H132").Value = 2441.1282
$ws.Range("I132").Value = 2104.1
$ws.Range("J132").Value = 3564.5557
$ws.Range("K132").Value = 6312.299999999999
$ws.Range("L132").Value = 10693.6671
$ws.Range("M132").Value = -3782.299999999999
$ws.Range("N132").Value = -15753.6671

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (BSM)
$ws.Range("H94").Value = 48630.81
$ws.Range("I94").Value = 756.0625
$ws.Range("K94").Value = 756.0625
$ws.Range("M94").Value = -305.0625

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 4054.862
$ws.Range("I31").Value = 3019.35
$ws.Range("J31").Value = 6356
$ws.Range("K31").Value = 3019.35
$ws.Range("L31").Value = 6356
$ws.Range("M31").Value = -2724.35
$ws.Range("N31").Value = -6946

# Row 34 (CRP)
$ws.Range("H34").Value = 4054.862
$ws.Range("I34").Value = 3019.35
$ws.Range("J34").Value = 6356
$ws.Range("K34").Value = 3019.35
$ws.Range("L34").Value = 6356
$ws.Range("M34").Value = -2817.35
$ws.Range("N34").Value = -6760

# Row 134 (CRP)
$ws.Range("H134").Value = 1636.72
$ws.Range("I134").Value = 1543.6522
$ws.Range("K134").Value = 4630.9566
$ws.Range("M134").Value = -2095.9566

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (CUL)
$ws.Range("H2").Value = 1075.4286
$ws.Range("I2").Value = 2252
$ws.Range("J2").Value = 421.77777
$ws.Range("K2").Value = 13512
$ws.Range("L2").Value = 2530.66662
$ws.Range("M2").Value = -13399
$ws.Range("N2").Value = -2756.66662

# Row 5 (CUL)
$ws.Range("H5").Value = 1467.3778
$ws.Range("I5").Value = 1155.375
$ws.Range("J5").Value = 1823.9524
$ws.Range("K5").Value = 3466.125
$ws.Range("L5").Value = 5471.857199999999
$ws.Range("M5").Value = -3354.125
$ws.Range("N5").Value = -5695.857199999999

# Row 9 (CUL)
$ws.Range("H9").Value = 2000
$ws.Range("J9").Value = 2000
$ws.Range("L9").Value = 6000
$ws.Range("N9").Value = -6448

# Row 12 (CUL)
$ws.Range("H12").Value = 1288300.8
$ws.Range("I12").Value = 28
$ws.Range("J12").Value = 1932437.2
$ws.Range("K12").Value = 84
$ws.Range("L12").Value = 5797311.6
$ws.Range("M12").Value = 89
$ws.Range("N12").Value = -5797657.6

# Row 13 (CUL)
$ws.Range("H13").Value = 250
$ws.Range("I13").Value = 250
$ws.Range("K13").Value = 750
$ws.Range("M13").Value = -582

# Row 113 (CUL)
$ws.Range("H113").Value = 4520.355
$ws.Range("J113").Value = 4766.5864
$ws.Range("L113").Value = 14299.7592
$ws.Range("N113").Value = -18639.7592

# Row 135 (CUL)
$ws.Range("H135").Value = 1467.3778
$ws.Range("I135").Value = 1155.375
$ws.Range("J135").Value = 1823.9524
$ws.Range("K135").Value = 10398.375
$ws.Range("L135").Value = 16415.5716
$ws.Range("M135").Value = -7863.375
$ws.Range("N135").Value = -21485.5716

$ws = $wb.Worksheets.Item("GSM")
# Row 109 (GSM)
$ws.Range("H109").Value = 20285
$ws.Range("J109").Value = 20285
$ws.Range("L109").Value = 20285
$ws.Range("N109").Value = -22365

$ws = $wb.Worksheets.Item("LTW")
# Row 55 (LTW)
$ws.Range("H55").Value = 320.07693
$ws.Range("I55").Value = 386.1
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 386.1
$ws.Range("L55").Value = 100
$ws.Range("M55").Value = -213.1
$ws.Range("N55").Value = -446

$ws = $wb.Worksheets.Item("WVR")
# Row 5 (WVR)
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 18 (WVR)
$ws.Range("H18").Value = 50000
$ws.Range("I18").Value = 50000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 50000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -49827
$ws.Range("N18").ClearContents()

# Row 113 (WVR)
$ws.Range("H113").Value = 926.2
$ws.Range("I113").Value = 1166.3334
$ws.Range("J113").Value = 823.2857
$ws.Range("K113").Value = 3499.0002
$ws.Range("L113").Value = 2469.8571
$ws.Range("M113").Value = -1329.0002
$ws.Range("N113").Value = -6809.8571

# Row 126 (WVR)
$ws.Range("H126").Value = 12198.429
$ws.Range("I126").Value = 13982.333
$ws.Range("J126").Value = 1495
$ws.Range("K126").Value = 41946.999
$ws.Range("L126").Value = 4485
$ws.Range("M126").Value = -39476.999
$ws.Range("N126").Value = -9425
